$d = $word.ActiveDocument

# RF004: Removido atributo sobrenome (remove the "Sobrenome;" list item paragraph)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Sobrenome;`r") {
        $p.Range.Delete()
        break
    }
}

# RF008: Removido máscara do atributo sobrenome (remove the "Sobrenome: Permitir apenas letras;" list item paragraph)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Sobrenome: Permitir apenas letras;`r") {
        $p.Range.Delete()
        break
    }
}

# RF015: Adicionado atributo cep (add "cep, " before "estado" in the Local list item)
$d.Content.Find.Execute("Local: estado, cidade, bairro, rua, número (opcional), complemento (opcional);", $true, $false, $false, $false, $false, $true, 1, $false, "Local: cep, estado, cidade, bairro, rua, número (opcional), complemento (opcional);", 2)
